# Auto-generated Excel COM-interop script
# Adds column-A metric values (2-factor full results) to all three sheets,
# updates the active-cell selection on each sheet, and adds page setup
# (paper size / orientation) to the AvgProcessorUtil sheet.

$wb = $excel.ActiveWorkbook

$accuracyValues = @(48.26,48.26,48.26,50.13,50.13,50.13,47.72,47.72,47.72,42.66,42.66,42.66,48.26,48.26,48.26,50.13,50.13,50.13,47.72,47.72,47.72,42.66,42.66,42.66,48.26,48.26,48.26,50.13,50.13,50.13,47.72,47.72,47.72,42.66,42.66,42.66)
$processorValues = @(26.15,26.65,26.625,26,23.6,23.274999999999999,25.849,26.6,26.125,22.65,23.625,23.1,50.674999999999997,45.9,46,46.55,49.875,45.3,45.024000000000001,46.024999999999999,46.225000000000001,44.575000000000003,47.024999999999999,48.973999999999997,45.3,46.5,46.674999999999997,45.024999999999999,45.4,50.598999999999997,45.375,48.424999999999997,49.55,44.875,43.774999999999999,42.973999999999997)
$trainTimeValues = @(281707.75,280915.25,280129.75,252211.75,251046.25,254584.25,236019.5,244353,237415,255814.25,259425.25,247506.25,136667.25,138424.25,136837.5,125970.75,126322.75,125614.5,121406,122068.5,122024.5,129412.5,128744.75,135448.5,138225.5,135742.5,134782.25,126399.5,124517.5,124249.25,120889.25,121840.5,120589.75,134018.25,128454.25,139104.25)

function Fill-ColumnA {
    param($ws, $values)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $values[$i]
    }
}

$wsAccuracy = $wb.Worksheets.Item("AvgAccuracy")
Fill-ColumnA $wsAccuracy $accuracyValues
$wsAccuracy.Range("G33").Select()

$wsProcessor = $wb.Worksheets.Item("AvgProcessorUtil")
Fill-ColumnA $wsProcessor $processorValues
$wsProcessor.PageSetup.PaperSize = 9
$wsProcessor.PageSetup.Orientation = 1
$wsProcessor.Range("G33").Select()

$wsTrainTime = $wb.Worksheets.Item("AvgTrainTime")
Fill-ColumnA $wsTrainTime $trainTimeValues
$wsTrainTime.Range("G33").Select()
